$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Periodo Mora" column (E16:E24) so the periods now read in
# ascending order (2006 .. 2102) instead of descending (2102 .. 2006).
$ws.Range("E16").Value = "2006"
$ws.Range("E17").Value = "2007"
$ws.Range("E18").Value = "2008"
$ws.Range("E19").Value = "2009"
$ws.Range("E20").Value = "2010"
$ws.Range("E21").Value = "2011"
$ws.Range("E22").Value = "2012"
$ws.Range("E23").Value = "2101"
$ws.Range("E24").Value = "2102"

# The "Valor Mora" (F) value that was tied to period 2102 moves along with
# it: it used to sit on row 16 (when 2102 was first) and now belongs to
# row 24 (where 2102 now is), while the rest of the rows keep 35112.
$ws.Range("F16").Value = 35112
$ws.Range("F24").Value = 25749
